$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price cells whose new values would otherwise be
# auto-detected as numbers by Excel (e.g. "243.95").
$numericLooking = @("D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D18","D19","D20","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (Coin / Link / Price / Volume columns).
$ws.Range("D2").Value = '30.424.11'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '1.880.90'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '243.95'
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").Value = '0.4714'
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").Value = '0.2879'
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("D9").Value = '0.06454'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").Value = '22.09'
$ws.Range("E10").Value = '  +0.40%  '
$ws.Range("D11").Value = '0.07786'
$ws.Range("E11").Value = '  +0.77%  '
$ws.Range("D12").Value = '1.892.68'
$ws.Range("E12").Value = '  +1.32%  '
$ws.Range("D13").Value = '95.35'
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").Value = '0.7220'
$ws.Range("E14").Value = '  -2.48%  '
$ws.Range("D15").Value = '5.161'
$ws.Range("E15").Value = '  -0.18%  '
$ws.Range("D16").Value = '279.03'
$ws.Range("E16").Value = '  +1.66%  '
$ws.Range("D17").Value = '30.402.85'
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").Value = '13.00'
$ws.Range("E18").Value = '  -1.89%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("D20").Value = '0.000007445'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").Value = '2.134.75'
$ws.Range("E21").Value = '  +0.84%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '5.242'
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("D24").Value = '6.267'
$ws.Range("E24").Value = '  +1.57%  '
$ws.Range("D25").Value = '163.89'
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("D26").Value = '9.040'
$ws.Range("E26").Value = '  -1.57%  '
$ws.Range("D27").Value = '18.76'
$ws.Range("E27").Value = '  +0.51%  '
$ws.Range("D28").Value = '1.884'
$ws.Range("E28").Value = '  -0.83%  '
$ws.Range("D29").Value = '1.334'
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("D30").Value = '0.09635'
$ws.Range("E30").Value = '  -3.26%  '
$ws.Range("D31").Value = '1.468'
$ws.Range("E31").Value = '  -2.62%  '
$ws.Range("D32").Value = '4.248'
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("D33").Value = '4.122'
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("D34").Value = '0.04841'
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("D35").Value = '1.122'
$ws.Range("E35").Value = '  +0.30%  '
$ws.Range("D36").Value = '0.6897'
$ws.Range("E36").Value = '  -0.24%  '
$ws.Range("D37").Value = '2.712'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").Value = '0.01879'
$ws.Range("E38").Value = '  +1.96%  '
$ws.Range("D39").Value = '2.814'
$ws.Range("E39").Value = '  +2.12%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '6.236'
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '74.35'
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("D42").Value = '1.951'
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").Value = '0.4243'
$ws.Range("E43").Value = '  +2.06%  '
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").Value = '0.8258'
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("D46").Value = '100.94'
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("D47").Value = '9.622'
$ws.Range("E47").Value = '  +3.15%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '35.13'
$ws.Range("E48").Value = '  -0.48%  '
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").Value = '6.927'
$ws.Range("E49").Value = '  -0.62%  '
$ws.Range("D50").Value = '901.39'
$ws.Range("E50").Value = '  -1.48%  '
$ws.Range("D51").Value = '0.05740'
$ws.Range("E51").Value = '  +1.45%  '

# Restore the default (unstyled) cell style now that the values are
# locked in as text - keeps the saved XML free of spurious style refs.
foreach ($addr in $numericLooking) {
    $ws.Range($addr).Style = "Normal"
}
